$d = $word.ActiveDocument

# 1. "-1db WRT300N" -> append router port description
$d.Content.Find.Execute("-1db WRT300N", $false, $false, $false, $false, $false, $true, 1, $false, "-1db WRT300N(1db Internet Port,3db Ethernet Port)", 2)

# 2. "-1db Switch(2960-24-TT)" -> append switch port description
$d.Content.Find.Execute("-1db Switch(2960-24-TT)", $false, $false, $false, $false, $false, $true, 1, $false, "-1db Switch(2960-24-TT)(24db FastEthernet Port,3db GigabitEthernet port)", 2)

# 3. "-2db PC" -> append PC port description
$d.Content.Find.Execute("-2db PC", $false, $false, $false, $false, $false, $true, 1, $false, "-2db PC(1db FastEthernet Port)", 2)

# 4. "-1db SmartPhone" -> append wireless description
$d.Content.Find.Execute("-1db SmartPhone", $false, $false, $false, $false, $false, $true, 1, $false, "-1db SmartPhone(Wireless)", 2)

# Move the hidden "_GoBack" bookmark from the "DHCP-OSI 3" paragraph to the end
# of the "-1db SmartPhone(Wireless)" paragraph (collapsed range, right after the
# appended text, before the paragraph mark).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "-1db SmartPhone(Wireless)" + [char]13) {
        $target = $p
    }
}
$tr = $target.Range
$endPos = $tr.End - 1
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Done"
